# Update "想去人数" (want-to-go count) figures refreshed from the source site.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 389
$ws1.Range("F9").Value = 556
$ws1.Range("F12").Value = 171
$ws1.Range("F13").Value = 13529
$ws1.Range("F14").Value = 180
$ws1.Range("F17").Value = 5563
$ws1.Range("F18").Value = 5587

# Sheet "全部类型" (sheet4) — aggregated view with the same rows duplicated
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F22").Value = 389
$ws4.Range("F31").Value = 556
$ws4.Range("F34").Value = 171
$ws4.Range("F35").Value = 13529
$ws4.Range("F36").Value = 180
$ws4.Range("F40").Value = 5563
$ws4.Range("F41").Value = 5587
